$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.344.72'
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").Value = '1.844.06'
$ws.Range("E3").Value = '  -0.18%  '
$ws.Range("D4").Value = '''0.9971'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '''239.98'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.19%  '
$ws.Range("D6").Value = '''0.6273'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("D7").Value = '''0.9991'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '''0.07482'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.53%  '
$ws.Range("D9").Value = '''0.2898'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.11%  '
$ws.Range("D10").Value = '''24.47'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.92%  '
$ws.Range("D11").Value = '''0.07738'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.08%  '
$ws.Range("D12").Value = '1.843.76'
$ws.Range("E12").Value = '  -2.35%  '
$ws.Range("D13").Value = '''4.982'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.80%  '
$ws.Range("D14").Value = '''0.6797'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.17%  '
$ws.Range("E15").Value = '  -0.62%  '
$ws.Range("D16").Value = '''81.95'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.31%  '
$ws.Range("D17").Value = '''6.184'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.93%  '
$ws.Range("D18").Value = '29.346.92'
$ws.Range("E18").Value = '  -0.08%  '
$ws.Range("D19").Value = '''229.28'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.59%  '
$ws.Range("D20").Value = '''12.32'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.15%  '
$ws.Range("D21").Value = '''0.9989'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").Value = '''7.506'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.15%  '
$ws.Range("D23").Value = '''0.9991'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("D24").Value = '''158.19'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.30%  '
$ws.Range("D25").Value = '''8.425'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("D26").Value = '''0.1367'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.98%  '
$ws.Range("D27").Value = '''17.52'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.74%  '
$ws.Range("D28").Value = '''0.06540'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +16.48%  '
$ws.Range("D29").Value = '''1.417'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.10%  '
$ws.Range("D30").Value = '''1.479'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.85%  '
$ws.Range("D31").Value = '''4.106'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.14%  '
$ws.Range("D32").Value = '''4.088'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.27%  '
$ws.Range("D33").Value = '''1.825'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.29%  '
$ws.Range("D34").Value = '''1.140'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.72%  '
$ws.Range("D35").Value = '''0.6944'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.25%  '
$ws.Range("D36").Value = '''2.579'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.18%  '
$ws.Range("D37").Value = '1.264.67'
$ws.Range("E37").Value = '  +2.93%  '
$ws.Range("D38").Value = '''2.834'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.08%  '
$ws.Range("D39").Value = '''0.01836'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.22%  '
$ws.Range("D40").Value = '''6.800'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.93%  '
$ws.Range("D41").Value = '''0.9166'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.20%  '
$ws.Range("D42").Value = '''0.9987'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.08%  '
$ws.Range("D43").Value = '2.002.77'
$ws.Range("E43").Value = '  +1.17%  '
$ws.Range("D44").Value = '''101.34'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.10%  '
$ws.Range("D45").Value = '''66.10'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.96%  '
$ws.Range("D46").Value = '''1.735'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.93%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '''0.00000000118'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.57%  '
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").Value = '''7.068'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.82%  '
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").Value = '''0.1161'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.88%  '
$ws.Range("D50").Value = '''0.3950'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.99%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '''8.955'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.32%  '
